# qa tables and OperationDefinition
# Update the "Comments" column (H) text for the data-element-mapping table,
# adjust the sheet view (zoom/selection), and set explicit (best-fit-like)
# column widths, matching the authored revision of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Updated "Comments" (column H) values -----------------------------
$ws.Range("H2").Value  = "Provider or Payer assigned tracking control number"
$ws.Range("H3").Value  = "Choice of ""claim"" or ""preauthorization"""
$ws.Range("H4").Value  = "Payer ID"
$ws.Range("H5").Value  = "Payer endpoint where the attachments are submitted using the `$submit-operation"
$ws.Range("H6").Value  = "Organization of Provider who submitted claim/prior authorization"
$ws.Range("H7").Value  = "Provider who submitted claim/prior authorization"
$ws.Range("H8").Value  = "claim/prior authorization ID (Provider or Payer Assigned)"
$ws.Range("H9").Value  = "claim/prior authorization line item numbers"
$ws.Range("H10").Value = "LOINC attachment codes"
$ws.Range("H11").Value = "Deadline form submitting  attachments  to  Payer"
$ws.Range("H12").Value = "Date of Service for claim/prior authorization"
$ws.Range("H13").Value = "Payer assigned patient identifier"
$ws.Range("H14").Value = "Patient Demographic information for patient matching"
$ws.Range("H15").Value = "Provider assigned patient identifer only for prior authorizatons"
$ws.Range("H16").Value = "Patient Demographic information for patient matching"

# --- Sheet view: zoom + selected cell ----------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("D20").Select() | Out-Null

# --- Column widths (approximate best-fit widths for the new content) ---
$ws.Columns.Item(1).ColumnWidth = 2.5833333333333335
$ws.Columns.Item(2).ColumnWidth = 34.25
$ws.Columns.Item(3).ColumnWidth = 14.75
$ws.Columns.Item(4).ColumnWidth = 16.25
$ws.Columns.Item(5).ColumnWidth = 7.75
$ws.Columns.Item(6).ColumnWidth = 31.75
$ws.Columns.Item(7).ColumnWidth = 40.083333333333336
$ws.Columns.Item(8).ColumnWidth = 70.08333333333333
